$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-8 (Farmer_ID, Field_ID, State, District, Date, MDO_ID, Coordinates, Area) ---

# Row 2
$ws.Range("A2").Value = 195924
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "PUNJAB"
$ws.Range("D2").Value = "Sangrur"
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = "30.570212250344493-75.74285175651312,30.569905965745214-75.74303448200226,30.569882005646246-75.74362356215715,30.569514231818104-75.74375532567501,30.569369315845066-75.7434394955635,30.568719500722363-75.74383411556482,30.568748368656713-75.74429914355278,30.568599121343443-75.7443031668663,30.56885171579133-75.74639361351728,30.56786644938405-75.7463912665844,30.567483656395908-75.74661422520876,30.567459407015466-75.74749901890755,30.567811888484847-75.74753656983376,30.56785547947027-75.74750572443008,30.56894928924193-75.7475533336401,30.568774638469527-75.74642177671194,30.568878274255056-75.74642948806286,30.568968053355785-75.74735719710588,30.569207656040888-75.74736792594194,30.569465156712038-75.7473474740982,30.56949056029944-75.74755433946848,30.56999718727271-75.74753690510988,30.570018837939116-75.74827283620834,30.57164695420592-75.74829362332821,30.57162934529319-75.74692469090223,30.57164031478016-75.74549742043018,30.570686542118036-75.74549876153469,30.570681057320773-75.74481144547462,30.570305203631953-75.74478730559349"
$ws.Range("H2").Value = 30.88992760971416

# Row 3
$ws.Range("A3").Value = 195917
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "HARYANA"
$ws.Range("D3").Value = "Panipat"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "29.226596119121965-76.97419289499521,29.226707889959474-76.97216380387545,29.226707889959474-76.97216380387545,29.228222343523118-76.97216715663671,29.228287005744345-76.97630479931831,29.226388377188908-76.97692941874266"
$ws.Range("H3").Value = 19.97040285883758

# Row 4
$ws.Range("A4").Value = 195924
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "PUNJAB"
$ws.Range("D4").Value = "Sangrur"
$ws.Range("F4").Value = 141
$ws.Range("G4").Value = "30.561027910921187-75.74838481843472,30.561075835486864-75.74718721210957,30.561116542478885-75.74593160301447,30.561785751431326-75.74597451835871,30.5617990316313-75.74507363140583,30.560524990256354-75.74504479765892,30.560478797642087-75.74590511620045,30.560472157451976-75.74612472206354,30.56030384291616-75.74612438678741,30.5602940269651-75.74766028672457,30.5604173036902-75.74765760451555,30.56040229107632-75.74912544339895"
$ws.Range("H4").Value = 8.294462181977583

# Row 5
$ws.Range("A5").Value = 195921
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "PUNJAB"
$ws.Range("D5").Value = "Bathinda"
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = "30.245564408739-74.89460952579975,30.247195338320516-74.8946088552475,30.247449344984325-74.892469458282,30.245426252032658-74.89221565425395"
$ws.Range("H5").Value = 10.93759842907927

# Row 6
$ws.Range("A6").Value = 195921
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "PUNJAB"
$ws.Range("D6").Value = "Bathinda"
$ws.Range("F6").Value = 84
$ws.Range("G6").Value = "30.252706295557765-74.89442210644484,30.252433767046792-74.89994678646326,30.25040903899494-74.8997201398015,30.250688813657817-74.8972736299038,30.251573892644416-74.89719182252884,30.25175577307152-74.89481639117002"
$ws.Range("H6").Value = 20.05406191465647

# Row 7
$ws.Range("A7").Value = 195918
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "PUNJAB"
$ws.Range("D7").Value = "Tarn Taran"
$ws.Range("F7").Value = 124
$ws.Range("G7").Value = "31.304026410369165-74.67083685100079,31.304022399825403-74.67218969017267,31.305022452274617-74.67217493802309,31.3056692866186-74.67220511287451,31.30724395733098-74.67217963188887,31.30724882711044-74.67077549546957"
$ws.Range("H7").Value = 11.61850832827422

# Row 8
$ws.Range("A8").Value = 195920
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "PUNJAB"
$ws.Range("D8").Value = "Tarn Taran"
$ws.Range("F8").Value = 123
$ws.Range("G8").Value = "31.307241952127583-74.6612811461091,31.307376300659758-74.66097939759493,31.307010207620777-74.66068904846907,31.30666301948074-74.66058008372784,31.305877831022826-74.66076247394085,31.304723955920455-74.66139312833546,31.302650497372454-74.66228663921356,31.302586041242847-74.66440759599209,31.30342196220387-74.66450180858374,31.303465219049713-74.66580770909785,31.30557446751293-74.66596528887749,31.305646942547913-74.66451320797205,31.306159708989416-74.66446191072464,31.3061508287005-74.66423492878677,31.30668908728157-74.66419905424118,31.307726923645042-74.66418463736774,31.3077243455396-74.6618477627635,31.30716489499391-74.661780372262"
$ws.Range("H8").Value = 48.37050751439378

# --- Date column (E2:E8): convert from inline-string text dates to real numeric
#     Excel date/time serials, formatted with a custom datetime number format. ---

$ws.Range("E2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$dateFmt = $ws.Range("E2").NumberFormat
$ws.Range("E3:E8").NumberFormat = $dateFmt

$ws.Range("E2").Value = 45096
$ws.Range("E3").Value = 45096
$ws.Range("E4").Value = 45096
$ws.Range("E5").Value = 45096
$ws.Range("E6").Value = 45096
$ws.Range("E7").Value = 45096
$ws.Range("E8").Value = 45035

# --- Remove the last two source rows (old r=9, r=10): the data now only
#     spans down to row 8. ---
$ws.Rows(9).Delete()
$ws.Rows(9).Delete()
